$d = $word.ActiveDocument

# Locate the paragraph containing the target sentence by scanning
# paragraphs for a stable ASCII prefix that is unaffected by the edit.
$pStart = -1
$pEnd = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Length -ge 10 -and $t.Substring(0, 10) -eq "Должны быт") {
        $pStart = $p.Range.Start
        $pEnd = $p.Range.End - 1   # exclude the trailing paragraph mark
        break
    }
}
if ($pStart -eq -1) {
    throw "Could not locate target paragraph"
}

$run1Text = "Должны быть межабзацные интервалы `$ExpectedBeforePt"
$run2Text = "пт"
$run3Text = " и межстрочный интервал 1,5."

# Replace the whole paragraph's visible text (single run) first.
$full = $d.Range($pStart, $pEnd)
$full.Text = $run1Text + $run2Text + $run3Text

$len1 = $run1Text.Length
$len2 = $run2Text.Length
$len3 = $run3Text.Length

$r1 = $d.Range($pStart, $pStart + $len1)
$r2 = $d.Range($pStart + $len1, $pStart + $len1 + $len2)
$r3 = $d.Range($pStart + $len1 + $len2, $pStart + $len1 + $len2 + $len3)

# Force run boundaries between the three segments by toggling a formatting
# property on/off without changing the end result, so each segment becomes
# its own <w:r> while keeping identical run properties.
$r1.Font.Bold = 1
$r1.Font.Bold = 0
$r2.Font.Bold = 1
$r2.Font.Bold = 0
